$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 41500
$ws.Range("J3").Value = 41500
$ws.Range("L3").Value = 41500
$ws.Range("N3").Value = -41728

$ws.Range("H17").Value = 810.3125
$ws.Range("J17").Value = 711.7857
$ws.Range("L17").Value = 2135.3571
$ws.Range("N17").Value = -2471.3571

$ws.Range("H40").Value = 1266.6666
$ws.Range("J40").Value = 1400
$ws.Range("L40").Value = 1400
$ws.Range("N40").Value = -1750

$ws.Range("H102").Value = 41500
$ws.Range("J102").Value = 41500
$ws.Range("L102").Value = 41500
$ws.Range("N102").Value = -47990

$ws.Range("H112").Value = 25001378
$ws.Range("I112").Value = 250000460
$ws.Range("J112").Value = 1478.6945
$ws.Range("K112").Value = 750001380
$ws.Range("L112").Value = 4436.083500000001
$ws.Range("M112").Value = -750000272
$ws.Range("N112").Value = -6652.083500000001

$ws.Range("H113").Value = 4194
$ws.Range("I113").Value = 1495.6666
$ws.Range("J113").Value = 4772.2144
$ws.Range("K113").Value = 1495.6666
$ws.Range("L113").Value = 4772.2144
$ws.Range("M113").Value = 1758.3334
$ws.Range("N113").Value = -11280.2144

$ws.Range("H116").Value = 328702.78
$ws.Range("I116").Value = 836257.5
$ws.Range("J116").Value = 8141.8945
$ws.Range("K116").Value = 836257.5
$ws.Range("L116").Value = 8141.8945
$ws.Range("M116").Value = -832815.5
$ws.Range("N116").Value = -15025.8945

$ws.Range("H129").Value = 822.09
$ws.Range("J129").Value = 899.22095
$ws.Range("L129").Value = 2697.66285
$ws.Range("N129").Value = -12697.66285

$ws.Range("H137").Value = 1289071.2
$ws.Range("I137").Value = 2071799.8
$ws.Range("J137").Value = 3160.1428
$ws.Range("K137").Value = 6215399.4
$ws.Range("L137").Value = 9480.428400000001
$ws.Range("M137").Value = -6212849.4
$ws.Range("N137").Value = -14580.4284

$ws.Range("H138").Value = 5867.52
$ws.Range("I138").Value = 750.9
$ws.Range("J138").Value = 7146.675
$ws.Range("K138").Value = 2252.7
$ws.Range("L138").Value = 21440.025
$ws.Range("M138").Value = 2887.3
$ws.Range("N138").Value = -31720.025

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4122.0386
$ws.Range("I32").Value = 3193.6365
$ws.Range("J32").Value = 9228.25
$ws.Range("K32").Value = 3193.6365
$ws.Range("L32").Value = 9228.25
$ws.Range("M32").Value = -2906.6365
$ws.Range("N32").Value = -9802.25

$ws.Range("H61").Value = 4564.4
$ws.Range("I61").Value = 1187.4783
$ws.Range("J61").Value = 15660
$ws.Range("K61").Value = 1187.4783
$ws.Range("L61").Value = 15660
$ws.Range("M61").Value = -975.4783
$ws.Range("N61").Value = -16084

$ws.Range("H74").Value = 5112.522
$ws.Range("I74").Value = 5825.067
$ws.Range("J74").Value = 3776.5
$ws.Range("K74").Value = 5825.067
$ws.Range("L74").Value = 3776.5
$ws.Range("M74").Value = -4951.067
$ws.Range("N74").Value = -5524.5

$ws.Range("H77").Value = 5112.522
$ws.Range("I77").Value = 5825.067
$ws.Range("J77").Value = 3776.5
$ws.Range("K77").Value = 29125.335
$ws.Range("L77").Value = 18882.5
$ws.Range("M77").Value = -24757.335
$ws.Range("N77").Value = -27618.5

$ws.Range("H88").Value = 4765569
$ws.Range("I88").Value = 7410407
$ws.Range("J88").Value = 4860
$ws.Range("K88").Value = 7410407
$ws.Range("L88").Value = 4860
$ws.Range("M88").Value = -7410001
$ws.Range("N88").Value = -5672

$ws.Range("H91").Value = 4765569
$ws.Range("I91").Value = 7410407
$ws.Range("J91").Value = 4860
$ws.Range("K91").Value = 7410407
$ws.Range("L91").Value = 4860
$ws.Range("M91").Value = -7409003
$ws.Range("N91").Value = -7668

$ws.Range("H103").Value = 34744.8
$ws.Range("J103").Value = 34744.8
$ws.Range("L103").Value = 34744.8
$ws.Range("N103").Value = -37088.8

$ws.Range("H132").Value = 1831.4642
$ws.Range("I132").Value = 1177.75
$ws.Range("K132").Value = 3533.25
$ws.Range("M132").Value = -1003.25

$ws.Range("H136").Value = 4564.4
$ws.Range("I136").Value = 1187.4783
$ws.Range("J136").Value = 15660
$ws.Range("K136").Value = 3562.4349
$ws.Range("L136").Value = 46980
$ws.Range("M136").Value = -1012.4349
$ws.Range("N136").Value = -52080

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1971.8572
$ws.Range("I86").Value = 1748.1111
$ws.Range("J86").Value = 2139.6667
$ws.Range("K86").Value = 1748.1111
$ws.Range("L86").Value = 2139.6667
$ws.Range("M86").Value = -625.1111000000001
$ws.Range("N86").Value = -4385.6667

$ws.Range("H89").Value = 1971.8572
$ws.Range("I89").Value = 1748.1111
$ws.Range("J89").Value = 2139.6667
$ws.Range("K89").Value = 8740.5555
$ws.Range("L89").Value = 10698.3335
$ws.Range("M89").Value = -3124.5555
$ws.Range("N89").Value = -21930.3335

$ws.Range("H94").Value = 1199.3182
$ws.Range("I94").Value = 971.8570999999999
$ws.Range("J94").Value = 1597.375
$ws.Range("K94").Value = 971.8570999999999
$ws.Range("L94").Value = 1597.375
$ws.Range("M94").Value = -520.8570999999999
$ws.Range("N94").Value = -2499.375

$ws.Range("H99").Value = 2254
$ws.Range("I99").Value = 1028
$ws.Range("K99").Value = 1028
$ws.Range("M99").Value = 470

$ws.Range("H134").Value = 3031.6912
$ws.Range("I134").Value = 1141.6296
$ws.Range("J134").Value = 10321.929
$ws.Range("K134").Value = 3424.8888
$ws.Range("L134").Value = 30965.787
$ws.Range("M134").Value = -889.8887999999997
$ws.Range("N134").Value = -36035.787

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2995
$ws.Range("I31").Value = 940.94446
$ws.Range("J31").Value = 7103.1113
$ws.Range("K31").Value = 940.94446
$ws.Range("L31").Value = 7103.1113
$ws.Range("M31").Value = -645.94446
$ws.Range("N31").Value = -7693.1113

$ws.Range("H34").Value = 2995
$ws.Range("I34").Value = 940.94446
$ws.Range("J34").Value = 7103.1113
$ws.Range("K34").Value = 940.94446
$ws.Range("L34").Value = 7103.1113
$ws.Range("M34").Value = -738.94446
$ws.Range("N34").Value = -7507.1113

$ws.Range("H106").Value = 34850
$ws.Range("J106").Value = 34850
$ws.Range("L106").Value = 34850
$ws.Range("N106").Value = -37374

$ws.Range("H107").Value = 588.1724
$ws.Range("J107").Value = 782.4545000000001
$ws.Range("L107").Value = 782.4545000000001
$ws.Range("N107").Value = -4622.4545

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 12978.556
$ws.Range("I39").Value = 8000
$ws.Range("J39").Value = 13600.875
$ws.Range("K39").Value = 24000
$ws.Range("L39").Value = 40802.625
$ws.Range("N39").Value = -41390.625
$ws.Range("M39").Value = -23706

$ws.Range("H41").Value = 1436.909
$ws.Range("I41").Value = 486.5
$ws.Range("J41").Value = 1980
$ws.Range("K41").Value = 1459.5
$ws.Range("L41").Value = 5940
$ws.Range("M41").Value = -1121.5
$ws.Range("N41").Value = -6616

$ws.Range("H58").Value = 1750
$ws.Range("J58").Value = 2000
$ws.Range("L58").Value = 6000
$ws.Range("N58").Value = -6256

$ws.Range("H113").Value = 525.9452
$ws.Range("I113").Value = 512.3269
$ws.Range("J113").Value = 559.6667
$ws.Range("K113").Value = 1536.9807
$ws.Range("L113").Value = 1679.0001
$ws.Range("M113").Value = 633.0192999999999
$ws.Range("N113").Value = -6019.0001

$ws.Range("H121").Value = 2421
$ws.Range("I121").Value = 220
$ws.Range("J121").Value = 2538.9106
$ws.Range("K121").Value = 660
$ws.Range("L121").Value = 7616.7318
$ws.Range("M121").Value = 650
$ws.Range("N121").Value = -10236.7318

$ws.Range("H132").Value = 2784.5417
$ws.Range("J132").Value = 3102.5715
$ws.Range("L132").Value = 27923.1435
$ws.Range("N132").Value = -32983.1435

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 30000
$ws.Range("J105").Value = 30000
$ws.Range("L105").Value = 30000
$ws.Range("N105").Value = -36988

$ws.Range("H132").Value = 3281.7058
$ws.Range("I132").Value = 1717
$ws.Range("J132").Value = 4377
$ws.Range("K132").Value = 5151
$ws.Range("L132").Value = 13131
$ws.Range("M132").Value = -2621
$ws.Range("N132").Value = -18191

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1444.591
$ws.Range("I46").Value = 1906.8334
$ws.Range("J46").Value = 1271.25
$ws.Range("K46").Value = 1906.8334
$ws.Range("L46").Value = 1271.25
$ws.Range("M46").Value = -1718.8334
$ws.Range("N46").Value = -1647.25

$ws.Range("H98").Value = 35000
$ws.Range("J98").Value = 35000
$ws.Range("L98").Value = 35000
$ws.Range("N98").Value = -40990

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H133").Value = 35320.25
$ws.Range("J133").Value = 35320.25
$ws.Range("L133").Value = 35320.25
$ws.Range("N133").Value = -40380.25

$ws.Range("H136").Value = 3570.6667
$ws.Range("I136").Value = 1716.619
$ws.Range("J136").Value = 6166.3335
$ws.Range("K136").Value = 5149.857
$ws.Range("L136").Value = 18499.0005
$ws.Range("M136").Value = -2599.857
$ws.Range("N136").Value = -23599.0005
